# Updated scheduler and GUI 5/19/20 4:58PM
#
# Adds a "Subject" column (D) to the Coords sheet, giving each building row
# its corresponding subject/department code abbreviation, switches the
# active sheet/tab from "Schedule" to "Coords", and sets the Coords sheet
# to print in portrait orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coords")

# New header cell D1 = "Subject", bold/size-14 font (matches the other
# header cells' font) but with no border (new cell style).
$ws.Range("D1").Value = "Subject"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").Font.Size = 14

# Subject code per building row (D2:D17), matching each row's Building (A).
$ws.Range("D2").Value  = "ECOM"   # Administration
$ws.Range("D3").Value  = "A"      # Albert O. Kuhn Library
$ws.Range("D4").Value  = "HIST"   # Arts & Humanities
$ws.Range("D5").Value  = "BIOL"   # Biological Sciences
$ws.Range("D6").Value  = "CMPE"   # Engineering
$ws.Range("D7").Value  = "ART"    # Fine Arts
$ws.Range("D8").Value  = "CMSC"   # Information Technology
$ws.Range("D9").Value  = "BTEC"   # Interdisciplinary Life S
$ws.Range("D10").Value = "MATH"   # Janet & Walter Sondheim
$ws.Range("D11").Value = "B"      # Lecture Hall 1
$ws.Range("D12").Value = "PYSC"   # Math & Psychology
$ws.Range("D13").Value = "CHEM"   # Meyerhoff Chemistry
$ws.Range("D14").Value = "PHYS"   # Physics
$ws.Range("D15").Value = "POLI"   # Public Policy
$ws.Range("D16").Value = "STAT"   # Sherman Hall
$ws.Range("D17").Value = "ENGL"   # University Center

# Set the sheet to portrait orientation for printing.
$ws.PageSetup.Orientation = 1

# Switch the active/selected tab from Schedule to Coords, with D17 selected.
$ws.Activate() | Out-Null
$ws.Range("D17").Select() | Out-Null
